# Trade #50 (overall trade #79, HighProbConvergence strategy) closes, and a
# new trade #108 (momentum strategy) opens.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet - roll up totals after the close of trade #79
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.49   # Current Capital
$summary.Range("B4").Value = 0.6       # Total P&L $
$summary.Range("B5").Value = 0.15      # Total P&L %
$summary.Range("B6").Value = 78        # Total Trades
$summary.Range("B7").Value = 39        # Winning Trades
$summary.Range("B9").Value = 50        # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - HighProbConvergence row (row 3)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C3").Value = 100.16     # Capital
$status.Range("D3").Value = 5          # Trades
$status.Range("E3").Value = 0.17       # P&L $
$status.Range("F3").Value = 0.16       # P&L %
$status.Range("G3").Value = 80         # Win Rate %

# ---------------------------------------------------------------------
# All Trades sheet - update the now-closed trade #79 (row 80)
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G80").Value = 0.67        # Exit Price
$allTrades.Range("H80").Value = "CLOSED"    # Status
$allTrades.Range("I80").Value = 7.4161      # P&L %
$allTrades.Range("J80").Value = 0.05        # P&L $
$allTrades.Range("K80").Value = 100.16      # Capital After
$allTrades.Range("L80").Value = "early_exit" # Exit Reason
$allTrades.Range("M80").Value = 0.12        # Duration (min)

# Append the newly opened trade #108 (momentum strategy) as row 109
$allTrades.Range("A109").Value = 108
$allTrades.Range("B109").Value = "'2026-02-18"
$allTrades.Range("C109").Value = "00:20:26"
$allTrades.Range("D109").Value = "momentum"
$allTrades.Range("E109").Value = "DOWN"
$allTrades.Range("F109").Value = 0.623742
$allTrades.Range("G109").Value = "'"
$allTrades.Range("H109").Value = "OPEN"
$allTrades.Range("I109").Value = 0
$allTrades.Range("J109").Value = 0
$allTrades.Range("K109").Value = 99.64873713109129
$allTrades.Range("L109").Value = "'"
$allTrades.Range("M109").Value = 0
$allTrades.Range("N109").Value = 0
$allTrades.Range("O109").Value = 0
$allTrades.Range("P109").Value = 0.9
$allTrades.Range("Q109").Value = "Downward momentum: -3.810% over 10 samples"

# ---------------------------------------------------------------------
# momentum sheet - append the newly opened trade #108 as row 27
# ---------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")
$momentum.Range("A27").Value = 108
$momentum.Range("B27").Value = "'2026-02-18"
$momentum.Range("C27").Value = "00:20:26"
$momentum.Range("D27").Value = "momentum"
$momentum.Range("E27").Value = "DOWN"
$momentum.Range("F27").Value = 0.623742
$momentum.Range("G27").Value = "'"
$momentum.Range("H27").Value = "OPEN"
$momentum.Range("I27").Value = 0
$momentum.Range("J27").Value = 0
$momentum.Range("K27").Value = 99.64873713109129
$momentum.Range("L27").Value = 0
$momentum.Range("M27").Value = 0
$momentum.Range("N27").Value = 0.9
$momentum.Range("O27").Value = "Downward momentum: -3.810% over 10 samples"
$momentum.Range("P27").Value = "'"
$momentum.Range("Q27").Value = 0

# ---------------------------------------------------------------------
# HighProbConvergence sheet - update the now-closed trade #79 (row 6)
# ---------------------------------------------------------------------
$hpc = $wb.Worksheets.Item("HighProbConvergence")
$hpc.Range("G6").Value = 0.67
$hpc.Range("H6").Value = "CLOSED"
$hpc.Range("I6").Value = 7.4161
$hpc.Range("J6").Value = 0.05
$hpc.Range("K6").Value = 100.16
$hpc.Range("P6").Value = "early_exit"
$hpc.Range("Q6").Value = 0.12
